# Update column F (dSF) values to reflect the re-pulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = -8
    7  = -2
    8  = -1
    9  = -2
    11 = 0
    13 = -2
    14 = 0
    21 = 1
    25 = -3
    27 = -1
    28 = -3
    29 = -1
    31 = -6
    32 = 2
    33 = -10
    34 = -6
    35 = -9
    38 = -6
    40 = -7
    42 = -2
    44 = -6
    45 = -3
    50 = 0
    54 = -4
    55 = 0
    57 = -2
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
